# Applies the "Trained reduced data models" edit to PerformanceData.xlsx
$wb = $excel.ActiveWorkbook

# --- T-Tests sheet: updated T.TEST results in H14/H15 (recalculated automatically
#     from the formulas once source data on "Class-Based Contrastive" sheet changes) ---

# --- Class-Based Contrastive sheet: fill in previously-empty Dice/Hausdorff values
#     for iterations 5 and 9 (rows 6 and 10) ---
$wsContrastive = $wb.Worksheets.Item("Class-Based Contrastive")
$wsContrastive.Range("F6").Value = 0.86250000000000004
$wsContrastive.Range("G6").Value = 16.856100000000001
$wsContrastive.Range("F10").Value = 0.83979999999999999
$wsContrastive.Range("G10").Value = 31.640699999999999

# Select the cell that was left selected after the edit on this sheet
$wsContrastive.Range("G11").Select()

# --- T-Tests sheet becomes the active/selected tab again ---
$wsTTests = $wb.Worksheets.Item("T-Tests")
$wsTTests.Activate()
$wsTTests.Range("H16").Select()
